# Automated map update (2025-09-27 09:16:30)
# Applies:
#   1) Updates to existing case in row 15 (Caso 4054).
#   2) Insertion of a brand-new case (Caso 5941) as a new row 33, shifting the
#      following rows down by one.
#   3) Append of a brand-new case (Caso 7325) as the new last row (83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update row 15 (Caso 4054 - AV AVELLANEDA 4020)
# ---------------------------------------------------------------------------
$ws.Range("E15").Value = "Pendiente ADM"
$ws.Range("G15").Value = "Pendiente"
$ws.Range("H15").Value = "Retirar columna ya traspasaron nodo"

# ---------------------------------------------------------------------------
# 2) Insert a new row at position 33 (Caso 5941 - COCHABAMBA 4090)
#    This pushes the previous rows 33..81 down to 34..82.
# ---------------------------------------------------------------------------
$ws.Rows("33:33").Insert()

$r33 = $ws.Range("A33:R33")
$r33.NumberFormat = "@"
$ws.Range("A33").Value = "5941"
$ws.Range("B33").Value = "5/26/2025"
$ws.Range("C33").Value = "COCHABAMBA 4090"
$ws.Range("D33").Value = "5"
$ws.Range("E33").Value = "806926861"
$ws.Range("F33").Value = "PEBCOM"
$ws.Range("G33").Value = "Pendiente"
$ws.Range("H33").Value = "Columna base podrida colocar r400 para pedir traspaso de fuente"
$ws.Range("I33").NumberFormat = "General"
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = "Cambio"
$ws.Range("K33").Value = "Fuente Teco"
$ws.Range("L33").Value = "Pasante"
$ws.Range("M33").NumberFormat = "General"
$ws.Range("M33").Value = -58.422268
$ws.Range("N33").NumberFormat = "General"
$ws.Range("N33").Value = -34.627754
$ws.Range("O33").Value = "Boedo"
$ws.Range("P33").Value = "Capital Sur"
$ws.Range("Q33").Value = "PPT-P"
$ws.Range("R33").Value = "Fuera de Poligono OVL"

# ---------------------------------------------------------------------------
# 3) Append a new row at the end, position 83 (Caso 7325 - SALAS 596)
# ---------------------------------------------------------------------------
$r83 = $ws.Range("A83:R83")
$r83.NumberFormat = "@"
$ws.Range("A83").Value = "7325"
$ws.Range("B83").Value = "9/26/2025"
$ws.Range("C83").Value = "SALAS 596"
$ws.Range("D83").Value = "7"
$ws.Range("E83").Value = "810015422"
$ws.Range("F83").Value = "PEBCOM"
$ws.Range("G83").Value = "Pendiente"
$ws.Range("H83").Value = "Cambio terminal"
$ws.Range("I83").NumberFormat = "General"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = "Cambio"
$ws.Range("K83").Value = "Sin equipos"
$ws.Range("L83").Value = "Terminal"
$ws.Range("M83").NumberFormat = "General"
$ws.Range("M83").Value = -58.434543
$ws.Range("N83").NumberFormat = "General"
$ws.Range("N83").Value = -34.632772
$ws.Range("O83").Value = "Boedo"
$ws.Range("P83").Value = "Capital Sur"
$ws.Range("Q83").Value = "PPT-S"
$ws.Range("R83").Value = "Fuera de Poligono OVL"
